# Correccion de errores (Angel Waidelich):
# "Matematica(Ortigoza Paul)" -> "Matematica (Ortigoza Paul)"
# (agrega un espacio antes del parentesis) en las celdas del horario de
# los Jueves donde dicta esa materia Ortigoza Paul.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Matematica (Ortigoza Paul) "
$ws.Range("E6").Value = "Matematica (Ortigoza Paul)"
$ws.Range("E7").Value = "Matematica (Ortigoza Paul) "

# Deja la celda E7 como activa/seleccionada, como quedo al guardar el archivo.
$ws.Range("E7").Select()
